$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking or percent-looking text that must
# remain stored as text (matching the original inline-string cell format).
# We briefly force a Text number format so Excel does not coerce the string into
# a number/percentage, then clear the format again so the cell keeps the default
# (unstyled) appearance from before the edit.
$textCells = @{
    'D2' = '278.30'
    'D3' = '27.23'
    'E3' = '1.95%'
    'D4' = '4.857'
    'E4' = '1.90%'
    'E5' = '1.92%'
    'D6' = '6.992'
    'E6' = '1.04%'
    'D7' = '1.199'
    'E7' = '-8.52%'
    'D8' = '0.8787'
    'E8' = '1.03%'
    'E9' = '-2.37%'
    'D10' = '0.05170'
    'E10' = '2.96%'
    'D11' = '0.07512'
    'E11' = '0.54%'
    'D12' = '0.02946'
    'E12' = '1.18%'
    'D13' = '0.08969'
    'E13' = '-0.88%'
    'D14' = '0.001567'
    'E14' = '0.17%'
    'D15' = '0.0006392'
    'E15' = '0.32%'
    'D16' = '0.006079'
    'E16' = '4.03%'
    'D17' = '3.477'
    'E17' = '0.78%'
    'D18' = '3.305'
    'E18' = '-0.20%'
    'D19' = '2.247'
    'E19' = '-1.60%'
    'D21' = '0.1324'
    'E21' = '0.40%'
    'E22' = '-0.46%'
    'D23' = '0.1505'
    'E23' = '9.01%'
    'D24' = '0.04397'
    'E24' = '0.29%'
    'D25' = '0.001176'
    'E25' = '0.66%'
    'D26' = '0.003893'
    'E26' = '5.00%'
    'E28' = '-1.69%'
    'D29' = '0.0001644'
    'E29' = '1.64%'
    'D40' = '0.04082'
    'E40' = '0.22%'
    'D41' = '0.006799'
    'E42' = '0.55%'
    'E43' = '-6.97%'
    'E44' = '0.25%'
    'D45' = '0.00005357'
    'E45' = '3.44%'
    'D46' = '1.561'
    'E46' = '5.01%'
    'D47' = '0.01852'
    'E47' = '-19.61%'
}

foreach ($addr in $textCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $textCells[$addr]
    $rng.ClearFormats()
}

# Plain text cells (coin names / links) can be set directly.
$ws.Range('B23').Value = 'ZBToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('B24').Value = 'CoinExToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
